# Applies the edits described by the commit diff:
#  1) Bold the "제품 설명" table header and append a (non-bold) ":" run.
#  2) Bold the "시장 점유율(%)" and "가격 책정 전략" table headers.
#  3) Ten body/table sentence rewrites (Korean localization tweaks).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Table 1, header row, column 2: "제품 설명" -> bold "제품 설명" + ":"
# ---------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$descHeader = $t1.Cell(1, 2)
$descStart = $descHeader.Range.Start
$descRange = $d.Range($descStart, $descStart + 5)
$descRange.Text = "제품 설명:"
$descBold = $d.Range($descStart, $descStart + 5)
$descBold.Font.Bold = 1

# ---------------------------------------------------------------------
# 2) Table 2, header row: bold "시장 점유율(%)" (col 2) and
#    "가격 책정 전략" (col 3)
# ---------------------------------------------------------------------
$t2 = $d.Tables.Item(2)

$shareHeader = $t2.Cell(1, 2)
$shareStart = $shareHeader.Range.Start
$shareRange = $d.Range($shareStart, $shareStart + 9)
$shareRange.Font.Bold = 1

$priceHeader = $t2.Cell(1, 3)
$priceStart = $priceHeader.Range.Start
$priceRange = $d.Range($priceStart, $priceStart + 8)
$priceRange.Font.Bold = 1

# ---------------------------------------------------------------------
# 3) Sentence-level text rewrites
# ---------------------------------------------------------------------
function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Exact `
    "정통 블렌드: 저희 차이는 프리미엄 홍차 잎과 계피, 카다몬, 정향, 생강, 후추 등 다양한 지상 향신료의 조화로운 조합입니다. " `
    "정통 블렌드: 저희 차이는 프리미엄 홍차 잎과 계피, 카다몬, 정향, 생강, 후추 등 다양한 시그니처 지상 향신료와의 조화로운 믹스로 이루어집니다. "

Replace-Exact `
    "건강 강화 성분: 신비 향신료 차이 차의 각 성분은 천연 건강상의 이점을 위해 선택됩니다. " `
    "건강 강화 성분: Mystic Spice 차이 티의 각 성분은 천연의 건강 혜택에 기반하여 엄선됩니다. "

Replace-Exact `
    "풍부한 아로마와 맛: 따뜻하고 매운 향기와 우리의 차이의 깊고 상쾌한 맛은 하루를 시작하거나 저녁에 긴장을 풀 수있는 완벽한 음료입니다. " `
    "풍부한 아로마와 맛: 저희 차이의 따뜻하고 매운 향기와 깊고 상쾌한 맛은 하루를 시작하거나 저녁에 긴장을 풀기에 완벽한 음료의 조건입니다. "

Replace-Exact `
    "다재다능한 양조 옵션: 차이 김이 뜨거워지거나, 상쾌한 아이스 티로, 크리미한 라떼를 좋아하든, 저희 블렌드는 모든 취향에 맞게 다재다능합니다. " `
    "다양한 브루잉 옵션: 김이 뜨겁게 올라오는 차이나, 상쾌한 아이스 티, 크리미한 라떼 등, 저희 블렌드는 모든 취향에 맞는 다양성을 제공합니다. "

Replace-Exact `
    "지속 가능한 공급: 지속 가능성을 위해 최선을 다하고 있으며, 우리는 유기농 농업을 실천하는 소규모 농장에서 재료를 공급하여 최고의 품질뿐만 아니라 지구의 복지를 보장합니다." `
    "지속 가능한 공급: 저희는 지속 가능성을 위해 최선을 다하고 있으며, 유기농 농업을 실천하는 소규모 농장에서 재료를 공급함으로써 최고의 품질뿐만 아니라 지구의 건강에도 이바지하고 있습니다."

Replace-Exact `
    "우아한 패키징: 신비주의 향신료 차이 차는 아름답게 디자인된 친환경 포장재로 제공되며, 차 애호가들에게 이상적인 선물이거나 호화로운 간식입니다." `
    "우아한 패키징: Mystic Spice 차이 티는 아름답게 디자인된 친환경 포장재로 제공되므로 차 애호가들에게 뿐만 아니라 자신에게도 이상적이고 고급스러운 선물이 될 수 있습니다."

Replace-Exact `
    "고객 만족 보장: Microsoft는 제품 뒤에 서서 만족도 보장을 제공합니다. " `
    "고객 만족 보장: 저희는 제품을 뒷받침하며 고객 만족을 보장합니다. "

Replace-Exact `
    "이상적인 대상: 차 애호가, 건강에 민감한 개인, 따뜻하고 매운 음료 애호가, 전통적인 인도 차이의 풍부한 맛을 탐구하고자하는 사람." `
    "이상적인 대상: 차 애호가, 건강에 민감한 사람, 따뜻하고 매운 음료 애호가, 전통적인 인도 차이의 풍부한 맛을 탐구하고자하는 모든 사람."

Replace-Exact `
    "Tetley: Tetley는 라틴 아메리카, 특히 시장 리더인 브라질에서 강력한 입지를 가진 영국 차 회사입니다. " `
    "Tetley: Tetley는 라틴 아메리카, 특히 시장 선두 주자인 브라질에서 강력한 입지를 가진 영국 차 회사입니다. "

Replace-Exact `
    "현지 브랜드: 마테 팩터, 차이 메이트, 차이 브라질 등 라틴 아메리카에서 차이 차 제품을 제공하는 여러 현지 브랜드도 있습니다. " `
    "현지 브랜드: Mate Factor, Chai Mate, Chai Brasil 등 라틴 아메리카에서 차이 티 제품을 제공하는 여러 현지 브랜드도 있습니다. "

Write-Host "done"
